# CHECKER DONE week13 end
# Fill in the week13 column (N) results for each checker row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# Extend column N to the same (best-fit) width as the preceding week columns
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Header date for week 13 (2015-11-30), formatted like the other week headers
$ws.Range("N1").Value = 42338
$ws.Range("N1").NumberFormat = "m/d/yy"

# Per-person results for week 13, centered like the other week columns
$ws.Range("N2").Value = 5
$ws.Range("N2").HorizontalAlignment = $xlCenter

$ws.Range("N3").Value = "5-"
$ws.Range("N3").HorizontalAlignment = $xlCenter

$ws.Range("N4").Value = 5
$ws.Range("N4").HorizontalAlignment = $xlCenter

$ws.Range("N5").Value = 4
$ws.Range("N5").HorizontalAlignment = $xlCenter

$ws.Range("N6").Value = "5-"
$ws.Range("N6").HorizontalAlignment = $xlCenter

$ws.Range("N7").Value = "4-"
$ws.Range("N7").HorizontalAlignment = $xlCenter

$ws.Range("N8").Value = "5-"
$ws.Range("N8").HorizontalAlignment = $xlCenter

$ws.Range("N9").Value = 5
$ws.Range("N9").HorizontalAlignment = $xlCenter

$ws.Range("N10").Value = 5
$ws.Range("N10").HorizontalAlignment = $xlCenter

$ws.Range("N11").Value = "-"
$ws.Range("N11").HorizontalAlignment = $xlCenter

# Move the active selection to N2, matching the edited workbook's cursor
$ws.Range("N2").Select()
